$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J32").Value = 5654
$ws.Range("L32").Value = 5654
$ws.Range("H32").Value = 9276.200000000001
$ws.Range("N32").Value = -6306
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H38").Value = 31.166666
$ws.Range("K38").Value = 93.49999800000001
$ws.Range("M38").Value = 278.500002
$ws.Range("I38").Value = 31.166666
$ws.Range("I40").Value = 7825
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H40").Value = 7825
$ws.Range("M40").Value = -7650
$ws.Range("K40").Value = 7825
$ws.Range("I58").Value = 348
$ws.Range("H58").Value = 2685.5334
$ws.Range("K58").Value = 1044
$ws.Range("M58").Value = -894
$ws.Range("L62").Value = 10099.8
$ws.Range("H62").Value = 7555.591
$ws.Range("N62").Value = -11347.8
$ws.Range("K62").Value = 6807.294
$ws.Range("M62").Value = -6183.294
$ws.Range("I62").Value = 6807.294
$ws.Range("J62").Value = 10099.8
$ws.Range("H65").Value = 7555.591
$ws.Range("K65").Value = 34036.47
$ws.Range("M65").Value = -30916.47
$ws.Range("N65").Value = -56739
$ws.Range("I65").Value = 6807.294
$ws.Range("J65").Value = 10099.8
$ws.Range("L65").Value = 50499
$ws.Range("J80").Value = 14966.923
$ws.Range("L80").Value = 44900.769
$ws.Range("H80").Value = 10962.167
$ws.Range("K80").Value = 1649.4
$ws.Range("N80").Value = -46896.769
$ws.Range("M80").Value = -651.3999999999999
$ws.Range("I80").Value = 549.8
$ws.Range("H83").Value = 10962.167
$ws.Range("K83").Value = 4948.2
$ws.Range("N83").Value = -144686.307
$ws.Range("M83").Value = 43.80000000000018
$ws.Range("I83").Value = 549.8
$ws.Range("J83").Value = 14966.923
$ws.Range("L83").Value = 134702.307
$ws.Range("J87").Value = 103451
$ws.Range("L87").Value = 103451
$ws.Range("H87").Value = 103451
$ws.Range("N87").Value = -105947
$ws.Range("K87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J90").Value = 103451
$ws.Range("L90").Value = 310353
$ws.Range("H90").Value = 103451
$ws.Range("K90").Value = 0
$ws.Range("N90").Value = -322833
$ws.Range("I90").Value = 0
$ws.Range("J116").Value = 16222
$ws.Range("L116").Value = 16222
$ws.Range("H116").Value = 1454529.1
$ws.Range("N116").Value = -23106
$ws.Range("H137").Value = 13958.267
$ws.Range("K137").Value = 59714.33099999999
$ws.Range("M137").Value = -57164.33099999999
$ws.Range("I137").Value = 19904.777
$ws.Range("L138").Value = 116422.758
$ws.Range("H138").Value = 20853.771
$ws.Range("N138").Value = -126702.758
$ws.Range("J138").Value = 38807.586
$ws.Range("N38").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("M87").ClearContents()
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 16162.514
$ws.Range("K32").Value = 16162.514
$ws.Range("H32").Value = 15405.3545
$ws.Range("M32").Value = -15875.514
$ws.Range("I45").Value = 2731.3
$ws.Range("J45").Value = 4952
$ws.Range("L45").Value = 4952
$ws.Range("H45").Value = 3841.65
$ws.Range("K45").Value = 2731.3
$ws.Range("M45").Value = -2354.3
$ws.Range("N45").Value = -5706
$ws.Range("H61").Value = 8811.941000000001
$ws.Range("N61").Value = -17566.857
$ws.Range("M61").Value = -2768.3
$ws.Range("K61").Value = 2980.3
$ws.Range("I61").Value = 2980.3
$ws.Range("J61").Value = 17142.857
$ws.Range("L61").Value = 17142.857
$ws.Range("N88").Value = -19014.334
$ws.Range("M88").Value = -1793.75
$ws.Range("K88").Value = 2199.75
$ws.Range("I88").Value = 2199.75
$ws.Range("J88").Value = 18202.334
$ws.Range("L88").Value = 18202.334
$ws.Range("H88").Value = 9058
$ws.Range("K91").Value = 2199.75
$ws.Range("N91").Value = -21010.334
$ws.Range("M91").Value = -795.75
$ws.Range("I91").Value = 2199.75
$ws.Range("J91").Value = 18202.334
$ws.Range("L91").Value = 18202.334
$ws.Range("H91").Value = 9058
$ws.Range("K122").Value = 3919.5
$ws.Range("H122").Value = 1480.4231
$ws.Range("M122").Value = -1469.5
$ws.Range("I122").Value = 1306.5
$ws.Range("H130").Value = 54710
$ws.Range("N130").Value = -64750
$ws.Range("J130").Value = 54710
$ws.Range("L130").Value = 54710
$ws.Range("K136").Value = 8940.900000000001
$ws.Range("N136").Value = -56528.571
$ws.Range("H136").Value = 8811.941000000001
$ws.Range("M136").Value = -6390.900000000001
$ws.Range("I136").Value = 2980.3
$ws.Range("J136").Value = 17142.857
$ws.Range("L136").Value = 51428.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2308.25
$ws.Range("J105").Value = 2998.3333
$ws.Range("L105").Value = 2998.3333
$ws.Range("N105").Value = -6492.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 1142
$ws.Range("I16").Value = 1142
$ws.Range("H16").Value = 1447.5151
$ws.Range("M16").Value = -855
$ws.Range("H31").Value = 2634134.2
$ws.Range("K31").Value = 4167662.5
$ws.Range("M31").Value = -4167367.5
$ws.Range("I31").Value = 4167662.5
$ws.Range("K34").Value = 4167662.5
$ws.Range("M34").Value = -4167460.5
$ws.Range("I34").Value = 4167662.5
$ws.Range("H34").Value = 2634134.2
$ws.Range("J58").Value = 2140.8572
$ws.Range("L58").Value = 2140.8572
$ws.Range("H58").Value = 1638.8286
$ws.Range("N58").Value = -2546.8572
$ws.Range("L62").Value = 8699.833000000001
$ws.Range("H62").Value = 8043.909
$ws.Range("N62").Value = -9947.833000000001
$ws.Range("K62").Value = 7797.9375
$ws.Range("M62").Value = -7173.9375
$ws.Range("I62").Value = 7797.9375
$ws.Range("J62").Value = 8699.833000000001
$ws.Range("H65").Value = 8043.909
$ws.Range("K65").Value = 38989.6875
$ws.Range("M65").Value = -35869.6875
$ws.Range("N65").Value = -49739.165
$ws.Range("I65").Value = 7797.9375
$ws.Range("J65").Value = 8699.833000000001
$ws.Range("L65").Value = 43499.165
$ws.Range("I86").Value = 77288.39999999999
$ws.Range("J86").Value = 24999.889
$ws.Range("L86").Value = 24999.889
$ws.Range("H86").Value = 52520.156
$ws.Range("N86").Value = -27245.889
$ws.Range("M86").Value = -76165.39999999999
$ws.Range("K86").Value = 77288.39999999999
$ws.Range("K89").Value = 386442
$ws.Range("N89").Value = -136231.445
$ws.Range("H89").Value = 52520.156
$ws.Range("M89").Value = -380826
$ws.Range("I89").Value = 77288.39999999999
$ws.Range("J89").Value = 24999.889
$ws.Range("L89").Value = 124999.445
$ws.Range("I113").Value = 1142
$ws.Range("H113").Value = 1447.5151
$ws.Range("K113").Value = 1142
$ws.Range("M113").Value = 1028
$ws.Range("J134").Value = 2233.3333
$ws.Range("L134").Value = 6699.999899999999
$ws.Range("H134").Value = 1372.091
$ws.Range("N134").Value = -11769.9999
$ws.Range("K134").Value = 3708.3156
$ws.Range("M134").Value = -1173.3156
$ws.Range("I134").Value = 1236.1052
$ws.Range("N136").Value = -11522.5716
$ws.Range("H136").Value = 1638.8286
$ws.Range("J136").Value = 2140.8572
$ws.Range("L136").Value = 6422.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J131").Value = 2083.2856
$ws.Range("L131").Value = 6249.8568
$ws.Range("H131").Value = 186192.34
$ws.Range("N131").Value = -16329.8568
$ws.Range("M131").Value = -795180.1799999999
$ws.Range("K131").Value = 800220.1799999999
$ws.Range("I131").Value = 266740.06

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J80").Value = 13902.571
$ws.Range("L80").Value = 13902.571
$ws.Range("H80").Value = 10301.637
$ws.Range("K80").Value = 4000
$ws.Range("N80").Value = -15898.571
$ws.Range("M80").Value = -3002
$ws.Range("I80").Value = 4000
$ws.Range("H83").Value = 10301.637
$ws.Range("K83").Value = 20000
$ws.Range("N83").Value = -79496.855
$ws.Range("M83").Value = -15008
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 13902.571
$ws.Range("L83").Value = 69512.855
$ws.Range("I113").Value = 4811
$ws.Range("H113").Value = 5973.6
$ws.Range("K113").Value = 4811
$ws.Range("M113").Value = -2641
$ws.Range("L122").Value = 10797.9999
$ws.Range("K122").Value = 6013.143
$ws.Range("H122").Value = 2203.75
$ws.Range("M122").Value = -3563.143
$ws.Range("N122").Value = -15697.9999
$ws.Range("I122").Value = 2004.381
$ws.Range("J122").Value = 3599.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 2535.6667
$ws.Range("H40").Value = 3083.4375
$ws.Range("M40").Value = -2399.6667
$ws.Range("K40").Value = 2535.6667
$ws.Range("H61").Value = 1750
$ws.Range("N61").Value = -4404
$ws.Range("M61").Value = -798
$ws.Range("K61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 4000
$ws.Range("L61").Value = 4000
$ws.Range("K93").Value = 3449.5
$ws.Range("M93").Value = -2201.5
$ws.Range("I93").Value = 3449.5
$ws.Range("H93").Value = 2957.4
$ws.Range("H100").Value = 3584.0715
$ws.Range("K100").Value = 3117.8
$ws.Range("M100").Value = -2576.8
$ws.Range("I100").Value = 3117.8
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("H113").Value = 1750
$ws.Range("K113").Value = 1000
$ws.Range("N113").Value = -8340
$ws.Range("M113").Value = 1170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 24855.334
$ws.Range("M81").Value = -68604.336
$ws.Range("K81").Value = 69665.336
$ws.Range("I81").Value = 34832.668
$ws.Range("I84").Value = 34832.668
$ws.Range("H84").Value = 24855.334
$ws.Range("M84").Value = -343022.68
$ws.Range("K84").Value = 348326.68
